$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 12825777
$ws.Range("I28").Value = 18525402
$ws.Range("J28").Value = 1618.375
$ws.Range("K28").Value = 18525402
$ws.Range("L28").Value = 1618.375
$ws.Range("M28").Value = -18524917
$ws.Range("N28").Value = -2588.375

$ws.Range("H70").Value = 1675.375
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1675.375
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 5026.125
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -5566.125

$ws.Range("H73").Value = 1675.375
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1675.375
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 5026.125
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -6898.125

$ws.Range("H135").Value = 2795.5
$ws.Range("I135").Value = 1048.1351
$ws.Range("J135").Value = 15726
$ws.Range("K135").Value = 9433.215899999999
$ws.Range("L135").Value = 141534
$ws.Range("M135").Value = -6898.215899999999
$ws.Range("N135").Value = -146604

$ws.Range("H137").Value = 786.9268
$ws.Range("I137").Value = 644.95
$ws.Range("J137").Value = 922.1429000000001
$ws.Range("K137").Value = 1934.85
$ws.Range("L137").Value = 2766.4287
$ws.Range("M137").Value = 615.1499999999999
$ws.Range("N137").Value = -7866.4287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 364812.25
$ws.Range("I32").Value = 2544.3699
$ws.Range("J32").Value = 4772404.5
$ws.Range("K32").Value = 2544.3699
$ws.Range("L32").Value = 4772404.5
$ws.Range("M32").Value = -2257.3699
$ws.Range("N32").Value = -4772978.5

$ws.Range("H61").Value = 1166.925
$ws.Range("I61").Value = 762.13043
$ws.Range("J61").Value = 1714.5883
$ws.Range("K61").Value = 762.13043
$ws.Range("L61").Value = 1714.5883
$ws.Range("M61").Value = -550.13043
$ws.Range("N61").Value = -2138.5883

$ws.Range("H74").Value = 987.0513
$ws.Range("I74").Value = 1031.4517
$ws.Range("J74").Value = 815
$ws.Range("K74").Value = 1031.4517
$ws.Range("L74").Value = 815
$ws.Range("M74").Value = -157.4517000000001
$ws.Range("N74").Value = -2563

$ws.Range("H77").Value = 987.0513
$ws.Range("I77").Value = 1031.4517
$ws.Range("J77").Value = 815
$ws.Range("K77").Value = 5157.2585
$ws.Range("L77").Value = 4075
$ws.Range("M77").Value = -789.2584999999999
$ws.Range("N77").Value = -12811

$ws.Range("H122").Value = 2016.9166
$ws.Range("I122").Value = 2197.3704
$ws.Range("J122").Value = 1475.5555
$ws.Range("K122").Value = 6592.111199999999
$ws.Range("L122").Value = 4426.666499999999
$ws.Range("M122").Value = -4142.111199999999
$ws.Range("N122").Value = -9326.666499999999

$ws.Range("H136").Value = 1166.925
$ws.Range("I136").Value = 762.13043
$ws.Range("J136").Value = 1714.5883
$ws.Range("K136").Value = 2286.39129
$ws.Range("L136").Value = 5143.7649
$ws.Range("M136").Value = 263.60871
$ws.Range("N136").Value = -10243.7649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7248997
$ws.Range("I31").Value = 7937936.5
$ws.Range("J31").Value = 15135.667
$ws.Range("K31").Value = 7937936.5
$ws.Range("L31").Value = 15135.667
$ws.Range("M31").Value = -7937641.5
$ws.Range("N31").Value = -15725.667

$ws.Range("H34").Value = 7248997
$ws.Range("I34").Value = 7937936.5
$ws.Range("J34").Value = 15135.667
$ws.Range("K34").Value = 7937936.5
$ws.Range("L34").Value = 15135.667
$ws.Range("M34").Value = -7937734.5
$ws.Range("N34").Value = -15539.667

$ws.Range("H58").Value = 893.26
$ws.Range("I58").Value = 689.5263
$ws.Range("J58").Value = 1538.4166
$ws.Range("K58").Value = 689.5263
$ws.Range("L58").Value = 1538.4166
$ws.Range("M58").Value = -486.5263
$ws.Range("N58").Value = -1944.4166

$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -34178

$ws.Range("H136").Value = 893.26
$ws.Range("I136").Value = 689.5263
$ws.Range("J136").Value = 1538.4166
$ws.Range("K136").Value = 2068.5789
$ws.Range("L136").Value = 4615.2498
$ws.Range("M136").Value = 481.4211
$ws.Range("N136").Value = -9715.2498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 629.5405
$ws.Range("I5").Value = 413.48386
$ws.Range("J5").Value = 1745.8334
$ws.Range("K5").Value = 1240.45158
$ws.Range("L5").Value = 5237.5002
$ws.Range("M5").Value = -1128.45158
$ws.Range("N5").Value = -5461.5002

$ws.Range("H23").Value = 395.23077
$ws.Range("I23").Value = 191.42857
$ws.Range("J23").Value = 633
$ws.Range("K23").Value = 574.28571
$ws.Range("L23").Value = 1899
$ws.Range("M23").Value = -339.28571
$ws.Range("N23").Value = -2369

$ws.Range("H41").Value = 500.2
$ws.Range("I41").Value = 534
$ws.Range("J41").Value = 485.7143
$ws.Range("K41").Value = 1602
$ws.Range("L41").Value = 1457.1429
$ws.Range("M41").Value = -1264
$ws.Range("N41").Value = -2133.1429

$ws.Range("H43").Value = 5460.8696
$ws.Range("I43").Value = 1866.6666
$ws.Range("J43").Value = 6000
$ws.Range("K43").Value = 5599.9998
$ws.Range("L43").Value = 18000
$ws.Range("M43").Value = -5485.9998
$ws.Range("N43").Value = -18228

$ws.Range("H112").Value = 444597.78
$ws.Range("I112").Value = 5593105.5
$ws.Range("J112").Value = 3297.1428
$ws.Range("K112").Value = 16779316.5
$ws.Range("L112").Value = 9891.428400000001
$ws.Range("M112").Value = -16778208.5
$ws.Range("N112").Value = -12107.4284

$ws.Range("H113").Value = 904.8816
$ws.Range("I113").Value = 575.58826
$ws.Range("J113").Value = 999.7627
$ws.Range("K113").Value = 1726.76478
$ws.Range("L113").Value = 2999.2881
$ws.Range("M113").Value = 443.23522
$ws.Range("N113").Value = -7339.2881

$ws.Range("H122").Value = 582.88
$ws.Range("J122").Value = 789.7857
$ws.Range("L122").Value = 7108.071300000001
$ws.Range("N122").Value = -12008.0713

$ws.Range("H131").Value = 5953205.5
$ws.Range("I131").Value = 721.8125
$ws.Range("J131").Value = 7353790
$ws.Range("K131").Value = 2165.4375
$ws.Range("L131").Value = 22061370
$ws.Range("M131").Value = 2874.5625
$ws.Range("N131").Value = -22071450

$ws.Range("H135").Value = 629.5405
$ws.Range("I135").Value = 413.48386
$ws.Range("J135").Value = 1745.8334
$ws.Range("K135").Value = 3721.35474
$ws.Range("L135").Value = 15712.5006
$ws.Range("M135").Value = -1186.35474
$ws.Range("N135").Value = -20782.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1883.6666
$ws.Range("I41").Value = 1325.5
$ws.Range("J41").Value = 3000
$ws.Range("K41").Value = 1325.5
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = -970.5
$ws.Range("N41").Value = -3710

$ws.Range("H43").Value = 12279.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12279.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 12279.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -12581.5

$ws.Range("H126").Value = 13890315
$ws.Range("I126").Value = 1262.4
$ws.Range("J126").Value = 23811068
$ws.Range("K126").Value = 3787.2
$ws.Range("L126").Value = 71433204
$ws.Range("M126").Value = -1317.2
$ws.Range("N126").Value = -71438144

$ws.Range("H132").Value = 1594.0233
$ws.Range("I132").Value = 1611.909
$ws.Range("J132").Value = 1535
$ws.Range("K132").Value = 4835.727000000001
$ws.Range("L132").Value = 4605
$ws.Range("M132").Value = -2305.727000000001
$ws.Range("N132").Value = -9665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 1451.25
$ws.Range("I18").Value = 1451.25
$ws.Range("K18").Value = 1451.25
$ws.Range("M18").Value = -1279.25

$ws.Range("H41").Value = 5099
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

$ws.Range("H136").Value = 4187.7856
$ws.Range("I136").Value = 3112.9
$ws.Range("J136").Value = 6875
$ws.Range("K136").Value = 9338.700000000001
$ws.Range("L136").Value = 20625
$ws.Range("M136").Value = -6788.700000000001
$ws.Range("N136").Value = -25725

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 5000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H62").Value = 4142.4736
$ws.Range("J62").Value = 4077.6667
$ws.Range("L62").Value = 4077.6667
$ws.Range("N62").Value = -5325.6667

$ws.Range("H65").Value = 4142.4736
$ws.Range("J65").Value = 4077.6667
$ws.Range("L65").Value = 20388.3335
$ws.Range("N65").Value = -26628.3335

$ws.Range("H126").Value = 2071.6667
$ws.Range("I126").Value = 1930.6666
$ws.Range("J126").Value = 2212.6667
$ws.Range("K126").Value = 5791.9998
$ws.Range("L126").Value = 6638.000100000001
$ws.Range("M126").Value = -3321.9998
$ws.Range("N126").Value = -11578.0001

$ws.Range("H132").Value = 16892698
$ws.Range("I132").Value = 19531990
$ws.Range("J132").Value = 1221
$ws.Range("K132").Value = 58595970
$ws.Range("L132").Value = 3663
$ws.Range("M132").Value = -58593440
$ws.Range("N132").Value = -8723
